{"js": "const replacements = [\n  [\"859\u00d79=\", \"229\u00d72=\"],\n  [\"283\u00d74=\", \"981\u00d77=\"],\n  [\"359\u00d77=\", \"730\u00d76=\"],\n  [\"604\u00d78=\", \"588\u00d79=\"],\n  [\"175\u00d74=\", \"324\u00d72=\"],\n  [\"430\u00d72=\", \"815\u00d75=\"],\n  [\"268\u00d79=\", \"920\u00d72=\"],\n  [\"819\u00d77=\", \"867\u00d72=\"],\n  [\"193\u00d73=\", \"909\u00d72=\"],\n  [\"938\u00d76=\", \"765\u00d77=\"],\n  [\"655\u00d77=\", \"654\u00d72=\"],\n  [\"912\u00d73=\", \"515\u00d73=\"],\n  [\"662\u00d75=\", \"238\u00d77=\"],\n  [\"688\u00d78=\", \"450\u00d72=\"],\n  [\"285\u00d72=\", \"805\u00d73=\"],\n  [\"791\u00d74=\", \"277\u00d75=\"],\n  [\"930\u00d77=\", \"291\u00d73=\"],\n  [\"644\u00d75=\", \"651\u00d73=\"],\n  [\"475\u00d73=\", \"116\u00d77=\"],\n  [\"761\u00d77=\", \"472\u00d72=\"],\n  [\"767\u00d79=\", \"628\u00d76=\"],\n  [\"299\u00d75=\", \"496\u00d79=\"],\n  [\"683\u00d75=\", \"227\u00d77=\"],\n  [\"680\u00d74=\", \"582\u00d76=\"],\n  [\"376\u00d74=\", \"421\u00d73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"859\u00d79=\"; New = \"229\u00d72=\" },\n    @{ Old = \"283\u00d74=\"; New = \"981\u00d77=\" },\n    @{ Old = \"359\u00d77=\"; New = \"730\u00d76=\" },\n    @{ Old = \"604\u00d78=\"; New = \"588\u00d79=\" },\n    @{ Old = \"175\u00d74=\"; New = \"324\u00d72=\" },\n    @{ Old = \"430\u00d72=\"; New = \"815\u00d75=\" },\n    @{ Old = \"268\u00d79=\"; New = \"920\u00d72=\" },\n    @{ Old = \"819\u00d77=\"; New = \"867\u00d72=\" },\n    @{ Old = \"193\u00d73=\"; New = \"909\u00d72=\" },\n    @{ Old = \"938\u00d76=\"; New = \"765\u00d77=\" },\n    @{ Old = \"655\u00d77=\"; New = \"654\u00d72=\" },\n    @{ Old = \"912\u00d73=\"; New = \"515\u00d73=\" },\n    @{ Old = \"662\u00d75=\"; New = \"238\u00d77=\" },\n    @{ Old = \"688\u00d78=\"; New = \"450\u00d72=\" },\n    @{ Old = \"285\u00d72=\"; New = \"805\u00d73=\" },\n    @{ Old = \"791\u00d74=\"; New = \"277\u00d75=\" },\n    @{ Old = \"930\u00d77=\"; New = \"291\u00d73=\" },\n    @{ Old = \"644\u00d75=\"; New = \"651\u00d73=\" },\n    @{ Old = \"475\u00d73=\"; New = \"116\u00d77=\" },\n    @{ Old = \"761\u00d77=\"; New = \"472\u00d72=\" },\n    @{ Old = \"767\u00d79=\"; New = \"628\u00d76=\" },\n    @{ Old = \"299\u00d75=\"; New = \"496\u00d79=\" },\n    @{ Old = \"683\u00d75=\"; New = \"227\u00d77=\" },\n    @{ Old = \"680\u00d74=\"; New = \"582\u00d76=\" },\n    @{ Old = \"376\u00d74=\"; New = \"421\u00d73=\" }\n)\n\nforeach ($rep in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $rep.Old\n    $find.Replacement.Text = $rep.New\n    $find.Execute([ref]$rep.Old, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$rep.New, 2)\n}\n"}
